$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ryzen 7")

# Insert a new row at position 8, shifting the existing rows 8-16 down to 9-17.
$ws.Rows("8").Insert() | Out-Null

# Row 7 keeps its request/second value (5920) but is relabeled from
# "C++ (*)" to "C++ dragon".
$ws.Range("A7").Value2 = "C++ dragon"

# The newly inserted row 8 becomes the "C++ asio CppServer" entry. Copy the
# date-format styling from the neighboring date cell first so D8 gets the
# same style id as the rest of column D.
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A8").Value2 = "C++ asio CppServer"
$ws.Range("B8").Value2 = 7084
$ws.Range("C8").Value2 = 10479
$ws.Range("D8").Value2 = 45933

# Refresh the sortState cache so it covers the new extent (A2:D14 / B2:B14).
$sort = $ws.Sort
$sort.SortFields.Clear() | Out-Null
$sort.SortFields.Add($ws.Range("B2:B14")) | Out-Null
$sort.SetRange($ws.Range("A2:D14")) | Out-Null
$sort.Header = 0
$sort.Apply() | Out-Null

# Point the chart series at the extended ranges.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(,'Ryzen 7'!`$A`$1:`$A`$14,'Ryzen 7'!`$B`$1:`$B`$14,1)"

# Grow the chart by one row's worth of height so its anchor extends from
# row 21 to row 22, matching the extra data row now behind it.
$rowHeight = $ws.Rows("8").RowHeight
$co.Height = $co.Height + $rowHeight

# Match the saved selection (D8 was the last-edited cell).
$ws.Range("D8").Select() | Out-Null

Write-Host "Done"
